$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G25 previously showed "nefunguje" with the unused red/orange-ish "fillId=5" style.
# It should now look like the other "DONE" cells (style copied from G3) and show "DONE".
$ws.Range("G3").Copy()
$ws.Range("G25").PasteSpecial(-4122)
$ws.Range("G25").Value2 = $ws.Range("G3").Value2

# B28 and B29 gain the "DONE" label (style already matches the DONE style, just was empty).
$ws.Range("B28").Value2 = $ws.Range("B30").Value2
$ws.Range("B29").Value2 = $ws.Range("B30").Value2

# B31 gains both the "DONE" style and the "DONE" label (previously empty with a different style).
$ws.Range("B30").Copy()
$ws.Range("B31").PasteSpecial(-4122)
$ws.Range("B31").Value2 = $ws.Range("B30").Value2

$excel.CutCopyMode = 0

# Move the active selection to B31 (was B30 before the edit).
$ws.Range("B31").Select()
